$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire rights for both movies.`n"
$ws.Range("D2").Value = "both_movies, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire rights for `"Oppenheimer`" has been successfully recorded.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded that both movies will be acquired.`n"
$ws.Range("D4").Value = "both_movies, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights to `"Oppenheimer`" have been successfully acquired for Friday's screening.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been made to acquire rights for `"Barbie.`"`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision process has concluded without an agreement on the movie for Friday.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision-making process did not result in a consensus on which movie to show, so no decision has been made regarding Friday's movie.`n"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected in this meeting.`n"
$ws.Range("D9").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded as no choice of a movie is possible without further discussion.`n"
$ws.Range("D10").Value = "no_decision, "
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: I have recorded the decision: there was no agreement on a movie to show on Friday.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision process did not lead to a selection for Friday's movie.`n"
$ws.Range("D13").Value = "no_decision, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision regarding which movie to acquire has not been made.`n"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision-making process concluded without any choice of a movie for Friday, leading to the conclusion of no decision.`n"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday resulted in no conclusive agreement. Thus, the outcome is recorded as `"no decision.`"`n"
$ws.Range("D19").Value = "no_decision, "
$ws.Range("C20").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding the movie for Friday.`n"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected during the meeting.`n"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been selected for acquisition.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for the screening on Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision has been recorded as no choice of a movie is possible without further discussion.`n"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both `"Oppenheimer`" and `"Barbie.`"`n"
$ws.Range("D25").Value = "both_movies, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" has been recorded successfully.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selection was made.`n"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded to select `"Barbie`" as the movie for Friday.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to show `"Barbie`" on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday has resulted in no agreement being reached.`n"
$ws.Range("D30").Value = "no_decision, "
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday was not made, so no action will be taken.`n"
$ws.Range("D31").Value = "no_decision, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded as a `"no decision`" regarding which movie to show on Friday.`n"
$ws.Range("D32").Value = "no_decision, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie could not be made, resulting in no selection.`n"
$ws.Range("D33").Value = "no_decision, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the selection of a movie for Friday.`n"
$ws.Range("D34").Value = "no_decision, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire rights for both movies.`n"
$ws.Range("D35").Value = "both_movies, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision to show `"Barbie`" has been recorded.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie has been selected for Friday.`n"
$ws.Range("D37").Value = "no_decision, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has been concluded with no definitive choice made.`n"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The committee did not come to a decision about which movie to show on Friday, so no action will be taken regarding movie rights acquisition.`n"
$ws.Range("D39").Value = "no_decision, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie has been selected for Friday.`n"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision to select a movie for Friday was not made, resulting in no acquisition for the movie rights.`n"
$ws.Range("D41").Value = "no_decision, "
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie selection for Friday.`n"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie`" for Friday's screening.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be made, so I will call the no_decision function.`n"
$ws.Range("D44").Value = "no_decision, "
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision regarding which movie to acquire rights for on Friday could not be made.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision.`"`n"
$ws.Range("D47").Value = "no_decision, "
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision-making process ended without a consensus on Friday's movie, resulting in no decision being made.`n"
$ws.Range("D48").Value = "no_decision, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for both movies.`n"
$ws.Range("D50").Value = "both_movies, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision about Friday's movie was reached in this conversation.`n"
$ws.Range("D51").Value = "no_decision, "
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded with no selection made for the movie to be shown on Friday.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired for showing.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision for Friday's movie has resulted in no agreement.`n"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been selected for the assembly.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision on which movie to show on Friday could not be made, as there was no clear agreement among the committee members.`n"
$ws.Range("D57").Value = "no_decision, "
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached on the movie to be shown on Friday.`n"
$ws.Range("D58").Value = "no_decision, "
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" for the movie showing on Friday has been successfully recorded.`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie`" as the movie that will be shown on Friday.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: No decision was made regarding which movie to show on Friday.`n"
$ws.Range("D61").Value = "no_decision, "
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision-making process did not lead to a conclusion about Friday's movie, so no movie will be acquired at this time.`n"
$ws.Range("D62").Value = "no_decision, "
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision process concluded without an agreement on which movie to show, resulting in no decision being made.`n"
$ws.Range("D63").Value = "no_decision, "
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer`" as the movie to be shown on Friday.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: Based on the criteria provided, there was no definitive decision made regarding the movie for Friday. Therefore, I will record the decision as a no decision.`n"
$ws.Range("D65").Value = "no_decision, "
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision process has concluded without a selection for Friday's movie.`n"
$ws.Range("D66").Value = "no_decision, "
$ws.Range("C67").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both `"Oppenheimer`" and `"Barbie.`"`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no consensus was reached regarding the movie to be shown on Friday.`n"
$ws.Range("D68").Value = "no_decision, "

Write-Host "Applied all changes"
